$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.863.97"
$ws.Range("E2").Value = "  -4.25%  "
$ws.Range("D3").Value = "3.323.17"
$ws.Range("E3").Value = "  -6.10%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "183.26"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -9.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "534.09"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.607"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("D8").Value = "3.317.87"
$ws.Range("E8").Value = "  -5.96%  "
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.624"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -5.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.20"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -9.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.135"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -6.72%  "
$ws.Range("E13").Value = "  -2.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.16"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -7.64%  "
$ws.Range("D15").Value = "3.837.74"
$ws.Range("E15").Value = "  -6.76%  "
$ws.Range("D16").Value = "3.314.27"
$ws.Range("E16").Value = "  -6.56%  "
$ws.Range("E17").Value = "  -5.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.85"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.62%  "
$ws.Range("D19").Value = "64.556.30"
$ws.Range("E19").Value = "  -4.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.19"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -6.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.966"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -7.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "376.25"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.86"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.77%  "
$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.35"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -5.88%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.33"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.95"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.51%  "
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.70"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.68"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.50"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "29.17"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -6.21%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.83"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.84%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "645.08"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -6.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.38"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "60.12"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -6.18%  "
$ws.Range("E36").Value = "  -5.89%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.397"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.10"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.58%  "
$ws.Range("D40").Value = "0.0₃0736"
$ws.Range("E40").Value = "  +6.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.995"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("E42").Value = "  -3.53%  "
$ws.Range("D43").Value = "2.913.89"
$ws.Range("E43").Value = "  -6.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.53"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.73"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -10.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0406"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.93"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +10.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.67"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.62"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -8.94%  "
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("E51").Value = "  +0.38%  "
